$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185132741928101
$ws.Range("B1").Value = 2.028560876846313
$ws.Range("C1").Value = 6.299643516540527
$ws.Range("D1").Value = 2.307458162307739
$ws.Range("E1").Value = 1.197398900985718
